$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellD = $ws.Range("D2")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "38.830.03"
$cellD.Style = $styleD
$ws.Range("E2").Value = "  +1.19%  "

$cellD = $ws.Range("D3")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.105.42"
$cellD.Style = $styleD
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.01%  "

$cellD = $ws.Range("D5")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "227.67"
$cellD.Style = $styleD
$ws.Range("E5").Value = "  -0.15%  "

$cellD = $ws.Range("D6")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.618"
$cellD.Style = $styleD
$ws.Range("E6").Value = "  +1.15%  "

$ws.Range("E7").Value = "  +3.18%  "

$ws.Range("E8").Value = "  -0.06%  "

$cellD = $ws.Range("D9")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.390"
$cellD.Style = $styleD
$ws.Range("E9").Value = "  +2.31%  "

$ws.Range("E10").Value = "  +1.24%  "

$ws.Range("E11").Value = "  -0.60%  "

$cellD = $ws.Range("D12")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "15.79"
$cellD.Style = $styleD
$ws.Range("E12").Value = "  +6.48%  "

$cellD = $ws.Range("D13")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.417.31"
$cellD.Style = $styleD
$ws.Range("E13").Value = "  +1.23%  "

$cellD = $ws.Range("D14")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "22.08"
$cellD.Style = $styleD
$ws.Range("E14").Value = "  -1.51%  "

$ws.Range("E15").Value = "  +3.15%  "

$cellD = $ws.Range("D16")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "5.53"
$cellD.Style = $styleD
$ws.Range("E16").Value = "  +1.79%  "

$cellD = $ws.Range("D17")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.094.67"
$cellD.Style = $styleD
$ws.Range("E17").Value = "  +0.52%  "

$cellD = $ws.Range("D18")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "38.857.84"
$cellD.Style = $styleD
$ws.Range("E18").Value = "  +1.33%  "

$cellD = $ws.Range("D19")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "71.75"
$cellD.Style = $styleD
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("E20").Value = "  +1.47%  "

$ws.Range("E21").Value = "  +1.87%  "

$cellD = $ws.Range("D22")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "228.31"
$cellD.Style = $styleD
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("E24").Value = "  -3.91%  "

$ws.Range("E25").Value = "  -0.14%  "

$cellD = $ws.Range("D26")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "9.95"
$cellD.Style = $styleD
$ws.Range("E26").Value = "  +5.91%  "

$cellD = $ws.Range("D27")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "171.89"
$cellD.Style = $styleD
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("E29").Value = "  +4.34%  "

$cellD = $ws.Range("D30")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "19.36"
$cellD.Style = $styleD
$ws.Range("E30").Value = "  +1.62%  "

$cellD = $ws.Range("D31")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.57"
$cellD.Style = $styleD
$ws.Range("E31").Value = "  +11.08%  "

$ws.Range("E32").Value = "  +0.90%  "

$ws.Range("E33").Value = "  +2.41%  "

$ws.Range("E34").Value = "  +13.58%  "

$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("E37").Value = "  +0.23%  "

$cellD = $ws.Range("D38")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "3.51"
$cellD.Style = $styleD
$ws.Range("E38").Value = "  -0.99%  "

$cellD = $ws.Range("D39")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "1.00"
$cellD.Style = $styleD
$ws.Range("E39").Value = "  +0.11%  "

$cellD = $ws.Range("D40")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "18.05"
$cellD.Style = $styleD
$ws.Range("E40").Value = "  -1.22%  "

$cellD = $ws.Range("D41")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "102.00"
$cellD.Style = $styleD
$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("E42").Value = "  +3.50%  "

$cellD = $ws.Range("D43")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "1.524.22"
$cellD.Style = $styleD
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("E44").Value = "  +8.27%  "

$cellD = $ws.Range("D46")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "7.82"
$cellD.Style = $styleD
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cellD = $ws.Range("D47")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.0918"
$cellD.Style = $styleD
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cellD = $ws.Range("D48")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "1.09"
$cellD.Style = $styleD
$ws.Range("E48").Value = "  +6.43%  "

$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("E50").Value = "  -0.09%  "

$cellD = $ws.Range("D51")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.305.42"
$cellD.Style = $styleD
$ws.Range("E51").Value = "  +1.27%  "
